# DESP_YR_FIN.xlsx update: insert a new "latest fiscal year" (FY2018,
# 2018-12-31) data column before column D on the DESP sheet, shifting all
# the existing yearly figures one column to the right, then populate the
# new column with the newly reported figures. Also corrects the
# "Capital Expenditures" (row 91) historical figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DESP")

# 1) Insert a new blank column at D; this shifts D:K -> E:L and keeps each
#    cell's value/type, but the brand-new column D cells default to the
#    general style, so formats get fixed up in step 2.
$ws.Columns("D").Insert()

# 2) Copy the number/date formatting from the (now shifted) column E onto
#    the new column D so the new cells render as dates/#,##0 numbers just
#    like the rest of their row, instead of "General".
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)

# 3) Populate the new column D with the newly reported FY2018 figures.
# -- Period Ending header rows --
$ws.Range("D7").Value = 43465
$ws.Range("D38").Value = 43465
$ws.Range("D80").Value = 43465

# -- Income Statement --
$ws.Range("D8").Value = 530600
$ws.Range("D9").Value = 172100
$ws.Range("D10").Value = 358500
$ws.Range("D12").Value = 71200
$ws.Range("D17").Value = 485200
$ws.Range("D18").Value = 45400
$ws.Range("D20").Value = -19200
$ws.Range("D21").Value = 41300
$ws.Range("D23").Value = 26200
$ws.Range("D24").Value = 7100
$ws.Range("D26").Value = 19200
$ws.Range("D27").Value = 19200
$ws.Range("D32").Value = 19200
$ws.Range("D33").Value = 19200
$ws.Range("D35").Value = 19200

# -- Balance Sheet --
$ws.Range("D41").Value = 346500
$ws.Range("D42").Value = "NA"
$ws.Range("D43").Value = 233700
$ws.Range("D45").Value = 74200
$ws.Range("D46").Value = 654300
$ws.Range("D48").Value = 19700
$ws.Range("D49").Value = 73700
$ws.Range("D52").Value = 12800
$ws.Range("D54").Value = 760500
$ws.Range("D57").Value = 227800
$ws.Range("D58").Value = 31200
$ws.Range("D59").Value = 126800
$ws.Range("D60").Value = 385700
$ws.Range("D62").Value = 127200
$ws.Range("D66").Value = 512900
$ws.Range("D72").Value = -306300
$ws.Range("D76").Value = 247600

# -- Cash Flow Statement --
$ws.Range("D81").Value = 19200
$ws.Range("D83").Value = 15100
$ws.Range("D89").Value = -17600

# "Capital Expenditures" row also got restated for the already-existing
# years, not just shifted.
$ws.Range("D91").Value = -13100
$ws.Range("E91").Value = -8700
$ws.Range("F91").Value = -4400
$ws.Range("G91").Value = -7100

$ws.Range("D94").Value = -26600
$ws.Range("D100").Value = -1300
$ws.Range("D101").Value = -13100
$ws.Range("D102").Value = -58600
